# Apply the "done counting for all regions" edit:
# The "Тогучинский" (Toguchin) region row had its file-name codes
# corrected from TOGUCHIN / sun_toguchin to TOGUCHI / sun_toguchi.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the row for the "Тогучинский" district in column A and update
# the corresponding meteo-data / sun-radiation file name codes in B/C.
$found = $false
$lastRow = $ws.Cells(1,1).End(4).Row  # xlDown = 4
for ($r = 1; $r -le $lastRow; $r++) {
    $districtCell = $ws.Cells.Item($r, 1)
    if ($districtCell.Value -eq "Тогучинский") {
        $ws.Cells.Item($r, 2).Value = "TOGUCHI"
        $ws.Cells.Item($r, 3).Value = "sun_toguchi"
        $found = $true
        break
    }
}

if (-not $found) {
    # Fallback: the district is known to live on row 24.
    $ws.Cells.Item(24, 2).Value = "TOGUCHI"
    $ws.Cells.Item(24, 3).Value = "sun_toguchi"
}

# Update the active selection to C24, matching the saved workbook state.
$ws.Range("C24").Select()
